# Enforce single work request per Excel file.
#
# This work report previously contained an extra "POL-40-2 / Pole,40ft,Class 2"
# line item that doesn't belong with this work request's billing period. The
# fix removes that line item, pulls the remaining line items up one row
# (content only -- the existing alternating row banding/style stays put), and
# refreshes the report's generated-on timestamp, summary totals and computed
# pricing to reflect the corrected set of line items. The row that used to
# hold the TOTAL (row 20) is no longer needed once the line items shift up, so
# the TOTAL now lives on row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates -------------------------------------------------

# Refresh the "Report Generated On" timestamp.
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Update report summary figures (total billed amount & line item count).
$ws.Range("C8").Value = 3096.6
$ws.Range("C9").Value = 3

# Scope ID # is no longer populated for this work request.
$ws.Range("G10").Value = ""

# --- Line item rows (content shifts up by one; formatting stays per-row) -----

# Row 16 <- old row 17 (PLA-DLOC / Inst / PLA,Difficult Location)
$ws.Range("A16").Value = "Point 11"
$ws.Range("B16").Value = "PLA-DLOC"
$ws.Range("C16").Value = "Inst"
$ws.Range("D16").Value = "PLA,Difficult Location"
$ws.Range("E16").Value = "EA"
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = 476.4

# Row 17 <- old row 18 (PLA-DLOC / Rem / PLA,Difficult Location)
$ws.Range("A17").Value = "Point 11"
$ws.Range("B17").Value = "PLA-DLOC"
$ws.Range("C17").Value = "Rem"
$ws.Range("D17").Value = "PLA,Difficult Location"
$ws.Range("E17").Value = "EA"
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = 476.4

# Row 18 <- old row 19 (PLA-BACK / Inst / Difficult Location Equip Adder-Backyard)
$ws.Range("A18").Value = "Point 11"
$ws.Range("B18").Value = "PLA-BACK"
$ws.Range("C18").Value = "Inst"
$ws.Range("D18").Value = "Difficult Location Equip Adder-Backyard"
$ws.Range("E18").Value = "EA"
$ws.Range("F18").Value = 18
$ws.Range("H18").Value = 2143.8

# --- TOTAL row moves from row 20 up to row 19 ---------------------------------

# Clear out the old line-item content that used to live on row 19 (now replaced
# by the TOTAL row), then clear the row that used to hold TOTAL (row 20).
$ws.Range("A19:H19").ClearContents()
$ws.Range("A20:H20").ClearContents()

# Re-merge A19:G19 for the TOTAL label (replaces the old A20:G20 merge) and
# write the new TOTAL row values.
$ws.Range("A20:G20").UnMerge()
$ws.Range("A19:G19").Merge()
$ws.Range("A19").Value = "TOTAL"
$ws.Range("H19").Value = 3096.6

# Update the worksheet's used-range dimension to match (A2:I19).
$ws.Range("A2:I19").Select()
